$wb = $excel.ActiveWorkbook

# The workbook has the same event data duplicated on two sheets:
#   "展览"   (Exhibitions)  - first sheet
#   "全部类型" (All types)   - fourth sheet
# Update the "想去人数" (interested-attendee count) column F for rows 3-5
# on both sheets to reflect the latest scraped numbers.

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 413
    $ws.Range("F4").Value = 25
    $ws.Range("F5").Value = 119
}
